# ============================================================
# Edit: blanko_renmin.docx
#  1) Insert a floating "PIS (Propam Integrated System)" text box
#     (printed-by stamp) as a new run before the very first run
#     of the document ("MARKAS BESAR").
#  2) Split the "NOMOR : ${no_surat}" table-cell paragraph so that
#     "${no_surat}" moves onto its own new paragraph.
# ============================================================

$d = $word.ActiveDocument

# ---------- Part 1: floating text box stamp ----------
$stampXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/>
          <w:b/>
          <w:bCs/>
          <w:noProof/>
          <w:kern w:val="0"/>
          <w:lang w:val="sv-SE"/>
        </w:rPr>
        <mc:AlternateContent>
          <mc:Choice Requires="wps">
            <w:drawing>
              <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251661312" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="5D47BB44" wp14:editId="082770BD">
                <wp:simplePos x="0" y="0"/>
                <wp:positionH relativeFrom="margin">
                  <wp:align>right</wp:align>
                </wp:positionH>
                <wp:positionV relativeFrom="paragraph">
                  <wp:posOffset>-114300</wp:posOffset>
                </wp:positionV>
                <wp:extent cx="2720340" cy="487680"/>
                <wp:effectExtent l="0" t="0" r="0" b="7620"/>
                <wp:wrapNone/>
                <wp:docPr id="718628983" name="Text Box 1"/>
                <wp:cNvGraphicFramePr/>
                <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
                  <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
                    <wps:wsp>
                      <wps:cNvSpPr txBox="1"/>
                      <wps:spPr>
                        <a:xfrm>
                          <a:off x="0" y="0"/>
                          <a:ext cx="2720340" cy="487680"/>
                        </a:xfrm>
                        <a:prstGeom prst="rect">
                          <a:avLst/>
                        </a:prstGeom>
                        <a:noFill/>
                        <a:ln w="6350">
                          <a:noFill/>
                        </a:ln>
                      </wps:spPr>
                      <wps:txbx>
                        <w:txbxContent>
                          <w:p>
                            <w:pPr>
                              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
                              <w:jc w:val="right"/>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                            </w:pPr>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>PIS (</w:t>
                            </w:r>
                            <w:proofErr w:type="spellStart"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>Propam</w:t>
                            </w:r>
                            <w:proofErr w:type="spellEnd"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t xml:space="preserve"> Integrated System)</w:t>
                            </w:r>
                          </w:p>
                          <w:p>
                            <w:pPr>
                              <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
                              <w:jc w:val="right"/>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                            </w:pPr>
                            <w:proofErr w:type="spellStart"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>Dicetak</w:t>
                            </w:r>
                            <w:proofErr w:type="spellEnd"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t xml:space="preserve"> </w:t>
                            </w:r>
                            <w:proofErr w:type="gramStart"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>Oleh :</w:t>
                            </w:r>
                            <w:proofErr w:type="gramEnd"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t xml:space="preserve"> ${user} | ${</w:t>
                            </w:r>
                            <w:proofErr w:type="spellStart"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>tgl_cetak</w:t>
                            </w:r>
                            <w:proofErr w:type="spellEnd"/>
                            <w:r>
                              <w:rPr>
                                <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                                <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                                <w:sz w:val="14"/>
                                <w:szCs w:val="14"/>
                                <w:lang w:val="en-GB"/>
                              </w:rPr>
                              <w:t>}</w:t>
                            </w:r>
                          </w:p>
                        </w:txbxContent>
                      </wps:txbx>
                      <wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="t" anchorCtr="0" forceAA="0" compatLnSpc="1">
                        <a:prstTxWarp prst="textNoShape">
                          <a:avLst/>
                        </a:prstTxWarp>
                        <a:noAutofit/>
                      </wps:bodyPr>
                    </wps:wsp>
                  </a:graphicData>
                </a:graphic>
                <wp14:sizeRelH relativeFrom="margin">
                  <wp14:pctWidth>0</wp14:pctWidth>
                </wp14:sizeRelH>
              </wp:anchor>
            </w:drawing>
          </mc:Choice>
          <mc:Fallback>
            <w:pict>
              <v:shapetype w14:anchorId="5D47BB44" id="_x0000_t202" coordsize="21600,21600" o:spt="202" path="m,l,21600r21600,l21600,xe">
                <v:stroke joinstyle="miter"/>
                <v:path gradientshapeok="t" o:connecttype="rect"/>
              </v:shapetype>
              <v:shape id="Text Box 1" o:spid="_x0000_s1026" type="#_x0000_t202" style="position:absolute;left:0;text-align:left;margin-left:163pt;margin-top:-9pt;width:214.2pt;height:38.4pt;z-index:251661312;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:right;mso-position-horizontal-relative:margin;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-width-relative:margin;v-text-anchor:top" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQAMg3A+FgIAACwEAAAOAAAAZHJzL2Uyb0RvYy54bWysU01vGyEQvVfqf0Dc6107juOsvI7cRK4q&#10;WUkkp8oZs+BdCRgK2Lvur+/A+ktpT1UvMPCG+XhvmD10WpG9cL4BU9LhIKdEGA5VY7Yl/fG2/DKl&#10;xAdmKqbAiJIehKcP88+fZq0txAhqUJVwBIMYX7S2pHUItsgyz2uhmR+AFQZBCU6zgEe3zSrHWoyu&#10;VTbK80nWgqusAy68x9unHqTzFF9KwcOLlF4EokqKtYW0urRu4brNZ6zYOmbrhh/LYP9QhWaNwaTn&#10;UE8sMLJzzR+hdMMdeJBhwEFnIGXDReoBuxnmH7pZ18yK1AuS4+2ZJv//wvLn/dq+OhK6r9ChgJGQ&#10;1vrC42Xsp5NOxx0rJYgjhYczbaILhOPl6G6U34wR4oiNp3eTaeI1u7y2zodvAjSJRkkdypLYYvuV&#10;D5gRXU8uMZmBZaNUkkYZ0pZ0cnObpwdnBF8ogw8vtUYrdJvu2MAGqgP25aCX3Fu+bDD5ivnwyhxq&#10;jPXi3IYXXKQCTAJHi5Ia3K+/3Ud/pB5RSlqcmZL6nzvmBCXqu0FR7ofjSENIh/Et0kKJu0Y214jZ&#10;6UfAsRziD7E8mdE/qJMpHeh3HO9FzIoQMxxzlzSczMfQTzJ+Dy4Wi+SEY2VZWJm15TF0pDNS+9a9&#10;M2eP/AdU7hlO08WKDzL0vr0Qi10A2SSNIsE9q0fecSSTdMfvE2f++py8Lp98/hsAAP//AwBQSwME&#10;FAAGAAgAAAAhAD14IJPfAAAABwEAAA8AAABkcnMvZG93bnJldi54bWxMj8FqwzAQRO+F/oPYQG+J&#10;HJMU4XodgiEUSntImktva2tjm1iSaymJ26+vempvO8ww8zbfTKYXVx595yzCcpGAYFs73dkG4fi+&#10;mysQPpDV1DvLCF/sYVPc3+WUaXeze74eQiNiifUZIbQhDJmUvm7ZkF+4gW30Tm40FKIcG6lHusVy&#10;08s0SR6loc7GhZYGLluuz4eLQXgpd2+0r1Kjvvvy+fW0HT6PH2vEh9m0fQIReAp/YfjFj+hQRKbK&#10;Xaz2okeIjwSE+VLFI9qrVK1AVAhrpUAWufzPX/wAAAD//wMAUEsBAi0AFAAGAAgAAAAhALaDOJL+&#10;AAAA4QEAABMAAAAAAAAAAAAAAAAAAAAAAFtDb250ZW50X1R5cGVzXS54bWxQSwECLQAUAAYACAAA&#10;ACEAOP0h/9YAAACUAQAACwAAAAAAAAAAAAAAAAAvAQAAX3JlbHMvLnJlbHNQSwECLQAUAAYACAAA&#10;ACEADINwPhYCAAAsBAAADgAAAAAAAAAAAAAAAAAuAgAAZHJzL2Uyb0RvYy54bWxQSwECLQAUAAYA&#10;CAAAACEAPXggk98AAAAHAQAADwAAAAAAAAAAAAAAAABwBAAAZHJzL2Rvd25yZXYueG1sUEsFBgAA&#10;AAAEAAQA8wAAAHwFAAAAAA==&#10;" filled="f" stroked="f" strokeweight=".5pt">
                <v:textbox>
                  <w:txbxContent>
                    <w:p>
                      <w:pPr>
                        <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
                        <w:jc w:val="right"/>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                      </w:pPr>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>PIS (</w:t>
                      </w:r>
                      <w:proofErr w:type="spellStart"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>Propam</w:t>
                      </w:r>
                      <w:proofErr w:type="spellEnd"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t xml:space="preserve"> Integrated System)</w:t>
                      </w:r>
                    </w:p>
                    <w:p>
                      <w:pPr>
                        <w:spacing w:before="0" w:after="0" w:line="240" w:lineRule="auto"/>
                        <w:jc w:val="right"/>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                      </w:pPr>
                      <w:proofErr w:type="spellStart"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>Dicetak</w:t>
                      </w:r>
                      <w:proofErr w:type="spellEnd"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t xml:space="preserve"> </w:t>
                      </w:r>
                      <w:proofErr w:type="gramStart"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>Oleh :</w:t>
                      </w:r>
                      <w:proofErr w:type="gramEnd"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t xml:space="preserve"> ${user} | ${</w:t>
                      </w:r>
                      <w:proofErr w:type="spellStart"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>tgl_cetak</w:t>
                      </w:r>
                      <w:proofErr w:type="spellEnd"/>
                      <w:r>
                        <w:rPr>
                          <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
                          <w:color w:val="747474" w:themeColor="background2" w:themeShade="80"/>
                          <w:sz w:val="14"/>
                          <w:szCs w:val="14"/>
                          <w:lang w:val="en-GB"/>
                        </w:rPr>
                        <w:t>}</w:t>
                      </w:r>
                    </w:p>
                  </w:txbxContent>
                </v:textbox>
                <w10:wrap anchorx="margin"/>
              </v:shape>
            </w:pict>
          </mc:Fallback>
        </mc:AlternateContent>
      </w:r>

</w:p>
'@

$anchor = $d.Content.Duplicate
$anchor.Find.Execute('MARKAS BESAR') | Out-Null
$anchor.Collapse(1)
$anchor.InsertXML($stampXml)

# ---------- Part 2: split "${no_surat}" onto its own paragraph ----------
$noSuratXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:right="56"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>NOMOR         :</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="0" w:after="0" w:line="276" w:lineRule="auto"/><w:ind w:right="56"/><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:eastAsia="Calibri" w:hAnsi="Arial Narrow" w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="id-ID"/></w:rPr><w:t>${no_surat}</w:t></w:r></w:p>
'@

$target = $d.Content.Duplicate
$target.Find.Execute('NOMOR         : ${no_surat}') | Out-Null
$targetRange = $d.Range($target.Start, $target.End)
$targetRange.InsertXML($noSuratXml)
